# Update the cryptos list worksheet with refreshed price / volume(1h) data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set D (Price) and E (Volume(1h)) for a given row as plain text.
# NOTE: positional parameters are used because named parameter binding
# (e.g. "-Row 2") is not reliable in this COM-interop PowerShell runtime.
# A leading apostrophe forces the Price cell to stay text even when the
# value looks numeric (e.g. "86.22"), matching the source inlineStr data.
function Set-Row($Row, $Price, $Volume) {
    if ($null -ne $Price) {
        $ws.Cells.Item($Row, 4).Value = "'" + $Price
    }
    if ($null -ne $Volume) {
        $ws.Cells.Item($Row, 5).Value = $Volume
    }
}

Set-Row 2  "40.424.95"  "  -3.12%  "
Set-Row 3  "2.359.86"   "  -4.54%  "
Set-Row 4  "0.999"      "  -0.01%  "
Set-Row 5  "308.87"     "  -3.05%  "
Set-Row 6  "86.22"      "  -7.33%  "
Set-Row 7  $null        "  -5.18%  "
Set-Row 9  "0.490"      "  -5.26%  "
Set-Row 10 "0.0834"     "  -4.83%  "
Set-Row 11 "30.33"      "  -8.27%  "
Set-Row 12 $null        "  -1.29%  "
Set-Row 13 "2.726.02"   "  -4.45%  "
Set-Row 14 $null        "  -6.21%  "
Set-Row 15 "14.88"      "  -4.47%  "
Set-Row 16 "2.375.75"   "  -2.81%  "
Set-Row 17 "0.752"      "  -6.20%  "
Set-Row 18 "40.342.58"  "  -3.21%  "
Set-Row 19 "0.0₃0902"   "  -5.05%  "
Set-Row 20 "6.07"       "  -6.46%  "
Set-Row 21 "68.14"      "  -4.20%  "
Set-Row 22 "10.63"      "  -6.16%  "
Set-Row 23 "232.74"     "  -3.42%  "
Set-Row 24 "2.57"       "  -6.86%  "
Set-Row 25 $null        "  +0.08%  "
Set-Row 26 $null        "  -8.51%  "
Set-Row 27 "23.39"      "  -6.68%  "
Set-Row 28 $null        "  -3.02%  "
Set-Row 29 "9.23"       "  -5.37%  "
Set-Row 30 "33.42"      "  -9.87%  "
Set-Row 31 "151.52"     "  -4.18%  "
Set-Row 32 $null        "  -0.01%  "
Set-Row 33 "5.15"       "  -6.54%  "
Set-Row 34 "2.43"       "  -4.93%  "
Set-Row 35 $null        "  -5.45%  "
Set-Row 36 $null        "  -2.75%  "

# Rows 37 and 38 swap coin identities (Celestia now ranks above LidoDAOToken).
$ws.Cells.Item(37, 2).Value = "Celestia"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(37, 4).Value = "'15.70"
$ws.Cells.Item(37, 5).Value = "  -9.42%  "

$ws.Cells.Item(38, 2).Value = "LidoDAOToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(38, 4).Value = "'2.72"
$ws.Cells.Item(38, 5).Value = "  -6.65%  "

Set-Row 39 "0.0984"     "  -5.92%  "
Set-Row 40 "1.69"       "  -9.53%  "
Set-Row 41 $null        "  -4.75%  "
Set-Row 42 $null        "  -5.70%  "
Set-Row 43 "1.946.70"   "  -2.49%  "
Set-Row 44 "0.0267"     "  -6.14%  "
Set-Row 45 $null        "  -8.79%  "
Set-Row 46 $null        "  -0.77%  "
Set-Row 47 "2.68"       "  -10.14%  "
Set-Row 48 "2.590.81"   "  -4.53%  "
Set-Row 49 "92.38"      "  -5.83%  "
Set-Row 50 "71.96"      "  -5.74%  "
Set-Row 51 "49.94"      "  -5.15%  "
